$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target table for rows 2..59 (A, B, C, D, E).
# A = client count index, B = auto scs (time in ms), C = auto capacity,
# D = dic_nbre_clients_poisson_2_keys, E = dic_nbre_clients_prob_poisson_2_values
$data = @(
    @(0, 33.94444444444444, 1.95, 0, 0.135),
    @(1, 33.94444444444444, 1.95, 3, 0.001),
    @(2, 33.94444444444444, 1.95, 4, 0.008),
    @(3, 33.94444444444444, 1.95, 5, 0.02),
    @(4, 33.94444444444444, 1.95, 6, 0.026),
    @(5, 33.94444444444444, 1.95, 7, 0.057),
    @(6, 33.94444444444444, 1.95, 8, 0.05),
    @(7, 33.94444444444444, 1.95, 9, 0.047),
    @(8, 33.94444444444444, 1.95, 10, 0.035),
    @(9, 33.94444444444444, 1.95, 11, 0.019),
    @(10, 33.94444444444444, 1.95, 12, 0.034),
    @(11, 33.94444444444444, 1.95, 13, 0.015),
    @(12, 33.94444444444444, 1.95, 14, 0.035),
    @(13, 33.94444444444444, 1.95, 15, 0.028),
    @(14, 33.94444444444444, 1.95, 16, 0.034),
    @(15, 33.94444444444444, 1.95, 17, 0.041),
    @(16, 33.94444444444444, 1.95, 18, 0.037),
    @(17, 33.94444444444444, 1.95, 19, 0.029),
    @(18, 33.94444444444444, 1.95, 20, 0.026),
    @(19, 33.94444444444444, 1.95, 21, 0.024),
    @(20, 33.94444444444444, 1.95, 22, 0.026),
    @(21, 33.94444444444444, 1.95, 23, 0.021),
    @(22, 33.94444444444444, 1.95, 24, 0.014),
    @(23, 33.94444444444444, 1.95, 25, 0.026),
    @(24, 33.94444444444444, 1.95, 26, 0.022),
    @(25, 33.94444444444444, 1.95, 27, 0.017),
    @(26, 33.94444444444444, 1.95, 28, 0.014),
    @(27, 33.94444444444444, 1.95, 29, 0.01),
    @(28, 33.94444444444444, 1.95, 30, 0.016),
    @(29, 33.94444444444444, 1.95, 31, 0.015),
    @(30, 33.94444444444444, 1.95, 32, 0.011),
    @(31, 33.94444444444444, 1.95, 33, 0.011),
    @(32, 33.94444444444444, 1.95, 34, 0.01),
    @(33, 33.94444444444444, 1.95, 35, 0.004),
    @(34, 33.94444444444444, 1.95, 36, 0.01),
    @(35, 33.94444444444444, 1.95, 37, 0.01),
    @(36, 33.94444444444444, 1.95, 38, 0.003),
    @(37, 33.94444444444444, 1.95, 39, 0.005),
    @(38, 33.94444444444444, 1.95, 40, 0.007),
    @(39, 33.94444444444444, 1.95, 41, 0.004),
    @(40, 33.94444444444444, 1.95, 42, 0.006),
    @(41, 33.94444444444444, 1.95, 43, 0.005),
    @(42, 33.94444444444444, 1.95, 44, 0.003),
    @(43, 33.94444444444444, 1.95, 45, 0.003),
    @(44, 33.94444444444444, 1.95, 46, 0.003),
    @(45, 33.94444444444444, 1.95, 47, 0.002),
    @(46, 33.94444444444444, 1.95, 48, 0.001),
    @(47, 33.94444444444444, 1.95, 49, 0.002),
    @(48, 33.94444444444444, 1.95, 50, 0.002),
    @(49, 33.94444444444444, 1.95, 51, 0.003),
    @(50, 33.94444444444444, 1.95, 52, 0.002),
    @(51, 33.94444444444444, 1.95, 53, 0.003),
    @(52, 33.94444444444444, 1.95, 54, 0.001),
    @(53, 33.94444444444444, 1.95, 57, 0.001),
    @(54, 33.94444444444444, 1.95, 58, 0.002),
    @(55, 33.94444444444444, 1.95, 59, 0.001),
    @(56, 33.94444444444444, 1.95, 74, 0.001),
    @(57, 33.94444444444444, 1.95, 77, 0.001)

)

$lastExistingRow = 54
$firstNewRow = 55

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $vals = $data[$i]

    if ($row -ge $firstNewRow) {
        # New row: clone the formatting of the last original data row (A54)
        # onto the new A-column cell so it keeps the bold/bordered/centered
        # style used by every other row in column A, then fill in values.
        $ws.Range("A" + $lastExistingRow).Copy() | Out-Null
        $ws.Range("A" + $row).PasteSpecial(-4122) | Out-Null
    }

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
}

$excel.CutCopyMode = 0
